$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header O1: rename "GI" -> "Carbon_nitrogen_ratio" and drop the wrap-text
# alignment that used to be applied to that header cell.
# ---------------------------------------------------------------------------
$ws.Range("O1").Value = "Carbon_nitrogen_ratio"
$ws.Range("O1").WrapText = $false

# ---------------------------------------------------------------------------
# O2:O8 data values were recomputed; also the column's number format reverts
# from a forced "0.00" back to the default General format.
# ---------------------------------------------------------------------------
$ws.Range("O2:O8").NumberFormat = "general"

$ws.Range("O2").Value = 24.82
$ws.Range("O3").Value = 23.94
$ws.Range("O4").Value = 23.07
$ws.Range("O5").Value = 22.1
$ws.Range("O6").Value = 21.51
$ws.Range("O7").Value = 19.54
$ws.Range("O8").Value = 13.1

# ---------------------------------------------------------------------------
# The active selection moved from P11 to K5.
# ---------------------------------------------------------------------------
$null = $ws.Range("K5").Select()
